$wb = $excel.ActiveWorkbook

# Sheet "VENTAS POR GRUPO"
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsGrupo.Range("M3").Value = 87.59
$wsGrupo.Range("M12").Value = "1 de 10"

# Sheet "VENTA MENSUAL"
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsMensual.Range("F3").Value = 594.47
$wsMensual.Range("F12").Value = 594.47
